# BDTData.xlsx -- "Major updates to analysis"
# Adds three new summary rows (19-21) describing the "Kaon5" / GTK extra-hits
# studies, tidies up the leftover helper formulas in rows 22-24, adds a new
# blank row 25, (re)creates a few merges, and sets page setup + selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel alignment constants (xlCenter = -4108)
$xlCenter = -4108

function Set-CenterVMiddle($rng) {
    $rng.HorizontalAlignment = $xlCenter
    $rng.VerticalAlignment = $xlCenter
}

# ---------------------------------------------------------------------
# Enter all the new text/number values first, in the same order the
# original author typed them, so shared-string indices line up.
# ---------------------------------------------------------------------
$ws.Range("D19").Value = "Kaon5"
$ws.Range("E19").Value = "^+ExtraHitsWord"
$ws.Range("F19").Value = 0.995
$ws.Range("G19").Value = 1807
$ws.Range("H19").Value = 31
$ws.Range("I19").Value = 18200
$ws.Range("J19").Value = 2974
$ws.Range("K19").Value = 8439
$ws.Range("L19").Value = 0.90069999999999995

$ws.Range("D20").Value = "Kaon5"
$ws.Range("E20").Value = "Same w/o hits veto applied such that BDT absorbs the veto"
$ws.Range("G20").Value = 1828
$ws.Range("H20").Value = 31
$ws.Range("I20").Value = 18399
$ws.Range("J20").Value = 2978
$ws.Range("K20").Value = 8533
$ws.Range("L20").Value = 0.90069999999999995

$ws.Range("D21").Value = "TMVA No GTK Extra Hits"
$ws.Range("G21").Value = 2366
$ws.Range("H21").Value = 60
$ws.Range("I21").Value = 18942
$ws.Range("J21").Value = 3165
$ws.Range("K21").Value = 8653
$ws.Range("L21").Value = 0.83460000000000001

$ws.Range("M19").Value = "One of these two rows HAS to be wrong -- sampleC, maybe I typed the stuff the wrong way"

$ws.Range("P20").Value = "Ke4 not tested for Kaon5"

# ---------------------------------------------------------------------
# Row 19 formatting
# ---------------------------------------------------------------------
$ws.Range("D19").Font.Color = 255
Set-CenterVMiddle $ws.Range("D19")
Set-CenterVMiddle $ws.Range("E19")
Set-CenterVMiddle $ws.Range("F19")
$ws.Range("G19").Font.Italic = $true
Set-CenterVMiddle $ws.Range("G19")
$ws.Range("H19").Font.Italic = $true
Set-CenterVMiddle $ws.Range("H19")
Set-CenterVMiddle $ws.Range("I19")
Set-CenterVMiddle $ws.Range("J19")
Set-CenterVMiddle $ws.Range("K19")
Set-CenterVMiddle $ws.Range("L19")
$ws.Range("M19:N20").WrapText = $true
Set-CenterVMiddle $ws.Range("M19:N20")
Set-CenterVMiddle $ws.Range("O19")

# ---------------------------------------------------------------------
# Row 20 formatting (taller row)
# ---------------------------------------------------------------------
$ws.Rows("20:20").RowHeight = 37

$ws.Range("D20").Font.Bold = $true
Set-CenterVMiddle $ws.Range("D20")
$ws.Range("E20:F20").WrapText = $true
Set-CenterVMiddle $ws.Range("E20:F20")
$ws.Range("G20").Font.Italic = $true
Set-CenterVMiddle $ws.Range("G20")
$ws.Range("H20").Font.Italic = $true
Set-CenterVMiddle $ws.Range("H20")
Set-CenterVMiddle $ws.Range("I20")
Set-CenterVMiddle $ws.Range("J20")
Set-CenterVMiddle $ws.Range("K20")
Set-CenterVMiddle $ws.Range("L20")
Set-CenterVMiddle $ws.Range("O20")

# ---------------------------------------------------------------------
# Row 21 formatting
# ---------------------------------------------------------------------
$ws.Range("D21").Font.Bold = $true
Set-CenterVMiddle $ws.Range("D21:F21")
$ws.Range("G21").Font.Italic = $true
Set-CenterVMiddle $ws.Range("G21")
$ws.Range("H21").Font.Italic = $true
Set-CenterVMiddle $ws.Range("H21")
Set-CenterVMiddle $ws.Range("I21")
Set-CenterVMiddle $ws.Range("J21")
Set-CenterVMiddle $ws.Range("K21")
Set-CenterVMiddle $ws.Range("L21")
Set-CenterVMiddle $ws.Range("M21")
Set-CenterVMiddle $ws.Range("N21")
Set-CenterVMiddle $ws.Range("O21")

# ---------------------------------------------------------------------
# Row 22 : drop the stray J22 cell entirely
# ---------------------------------------------------------------------
$ws.Range("J22").Clear()

# ---------------------------------------------------------------------
# Row 23 : the old %-difference helper formulas are no longer needed
# ---------------------------------------------------------------------
$ws.Range("I23:K23").ClearContents()

# ---------------------------------------------------------------------
# Row 24 : match vertical-centering used elsewhere in the block
# ---------------------------------------------------------------------
Set-CenterVMiddle $ws.Range("I24:K24")

# ---------------------------------------------------------------------
# Row 25 : new blank row under the block
# ---------------------------------------------------------------------
Set-CenterVMiddle $ws.Range("I25:K25")

# ---------------------------------------------------------------------
# Merges
# ---------------------------------------------------------------------
$ws.Range("D21:F21").Merge()
$ws.Range("M19:N20").Merge()
$ws.Range("E20:F20").Merge()

# ---------------------------------------------------------------------
# Page setup + selection
# ---------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("K31").Select()
